$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 488
$ws.Range("F6").Value = 886
$ws.Range("F9").Value = 2163
$ws.Range("F11").Value = 282
$ws.Range("F13").Value = 1059
$ws.Range("F14").Value = 177
$ws.Range("F15").Value = 2190
$ws.Range("F16").Value = 649
$ws.Range("F17").Value = 12384
$ws.Range("F18").Value = 1226
$ws.Range("F19").Value = 5
$ws.Range("F20").Value = 553
$ws.Range("F21").Value = 124
$ws.Range("F22").Value = 17
$ws.Range("F25").Value = 260
$ws.Range("F28").Value = 18
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 17
$ws.Range("F8").Value = 1
$ws.Range("F10").Value = 15
$ws.Range("F20").Value = 2
$ws.Range("F21").Value = 3
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5691
$ws.Range("F4").Value = 461
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5691
$ws.Range("F5").Value = 461
$ws.Range("F8").Value = 488
$ws.Range("F9").Value = 886
$ws.Range("F12").Value = 17
$ws.Range("F13").Value = 2163
$ws.Range("F15").Value = 282
$ws.Range("F19").Value = 1059
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = 177
$ws.Range("F23").Value = 15
$ws.Range("F24").Value = 2190
$ws.Range("F25").Value = 649
$ws.Range("F28").Value = 1226
$ws.Range("F29").Value = 5
$ws.Range("F30").Value = 553
$ws.Range("F31").Value = 124
$ws.Range("F32").Value = 17
$ws.Range("F38").Value = 260
$ws.Range("F43").Value = 2
$ws.Range("F45").Value = 3
$ws.Range("F49").Value = 18
